$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Akan/English parallel-text table had two story beats whose Akan and
# English lines were off-by-one-ish merged together. Insert a row to make
# room (18 -> 19 data rows) and rewrite every row so each Akan line lines up
# with its correct English line.
$ws.Rows("7:7").Insert()

$ws.Range("A1").Value = 'AKAN'
$ws.Range("B1").Value = 'ENGLISH'
$ws.Range("A2").Value = 'Ye'' nse se, nse se'
$ws.Range("B2").Value = 'We do not really mean, we do not really mean, (that what we are going to say is true)'
$ws.Range("A3").Value = 'SE YOYE A AKETEKYIRE, SE HYEYE'
$ws.Range("B3").Value = 'HOW AKETEKYIRE, THE CRICKET, GOT HIS TEETH BURNED '
$ws.Range("A4").Value = 'AKETEKYIRE ne Ananse na esiim'' se ye''ko pe aburo abedua. '
$ws.Range("B4").Value = 'AETEKYIRE, the Cricket, and Ananse, the Spider, started off, saying they were going in search of corn in order to plant it. '
$ws.Range("A5").Value = 'Na ye''nyaa aburo mmienu pe. '
$ws.Range("B5").Value = 'And they got just two grains of corn. '
$ws.Range("A6").Value = 'Na Ananse see, "Me''ko dua me die."'
$ws.Range("B6").Value = ' And Ananse said, " I am going off to plant mine." '
$ws.Range("A7").Value = 'Na Aketekyire see, "Me''ko kye me die m''awe." '
$ws.Range("B7").Value = 'And the Cricket said, " I am going to roast mine in order to chew it." '
$ws.Range("A8").Value = 'Ananse ko duaa ne die. '
$ws.Range("B8").Value = 'Ananse went and planted his. '
$ws.Range("A9").Value = 'Aketekyire so de ne die ko sii ogya so. '
$ws.Range("B9").Value = 'The Cricket, he, too, took his, and put it on the fire. '
$ws.Range("A10").Value = 'Emenee da, na osi ho.'
$ws.Range("B10").Value = 'It was long in getting cooked, and he sat down there (beside it). '
$ws.Range("A11").Value = 'Ananse so koo afuom''; osee o''fe ani afwe. '
$ws.Range("B11").Value = 'The Spider, too, he went off to the plantation; he was about to peep under the sheath covering the corn. '
$ws.Range("A12").Value = 'Aburo see, "Mfe m''ani, na wo nsa bewo m''ani." '
$ws.Range("B12").Value = 'The Corn said, " Don''t peep into my eye, lest your finger goes into my eye." '
$ws.Range("A13").Value = 'Ananse se o''twa n''ani afwe ne ''Kraman; '
$ws.Range("B13").Value = 'Ananse was turning his eyes to look at his Dog; '
$ws.Range("A14").Value = 'ne ''Kraman see, "Nfwe me, na nye me na me yee." '
$ws.Range("B14").Value = 'his Dog said, "You need not look at me, it was not I who did it." '
$ws.Range("A15").Value = 'Ananse se o''bu aba abefwe ''Kraman no, na Aba see, " Mmu me, na wo ''Kraman beka me." '
$ws.Range("B15").Value = 'Ananse was about to break off a stick to beat the Dog, when the stick said, "Don''t break me, for (if you do) your Dog will bite me." '
$ws.Range("A16").Value = 'Ananse baa ''fie beboo Aketekyire amanee. '
$ws.Range("B16").Value = 'Ananse came home and told the Cricket all about it. '
$ws.Range("A17").Value = 'Ose o''sere na oda ''gyam''. '
$ws.Range("B17").Value = 'While he was laughing, he fell into the fire. '
$ws.Range("A18").Value = 'Ene se Aketekyire ''se hyeye.'
$ws.Range("B18").Value = 'That is how the Cricket got his teeth burned.'
$ws.Range("A19").Value = 'M''anansesem a metooye yi, se eye de o, se ennye de o, momfa bi nko, na. momfa bi mmera.'
$ws.Range("B19").Value = 'This, my story, which I have related, if it be sweet, (or) if it be not sweet, take some elsewhere, and let some come back to me.'

# Rows with two-line wrapped text end up taller once the text is realigned.
$ws.Rows("11:11").RowHeight = 31
$ws.Rows("12:12").RowHeight = 14.5
$ws.Rows("15:15").RowHeight = 31
$ws.Rows("19:19").RowHeight = 31

# Narrower columns and a bigger zoom level, matching the re-saved workbook view.
$ws.Range("A:B").ColumnWidth = 95.54296875
$ws.Application.ActiveWindow.Zoom = 120
$ws.Range("A19").Select()
